# Update cryptocurrency price/volume figures (Price = column D, Volume(1h) = column E)
# for the rows whose scraped values changed in this run.
#
# Values are assigned via a temporary text-producing formula ( ="<text>" ) and then
# frozen in place with Copy + PasteSpecial(xlPasteValues). This guarantees the result
# is stored as plain text (matching the workbook's existing text cells) even for
# values that look like plain numbers (e.g. "526.77"), which a direct
# `.Value = "526.77"` assignment would otherwise auto-convert to a numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=""58.352.91"""
$ws.Range("E2").Formula = "=""  -3.16%  """
$ws.Range("D3").Formula = "=""3.141.84"""
$ws.Range("E3").Formula = "=""  -4.74%  """
$ws.Range("E4").Formula = "=""  +0.03%  """
$ws.Range("D5").Formula = "=""526.77"""
$ws.Range("E5").Formula = "=""  -5.40%  """
$ws.Range("D6").Formula = "=""133.76"""
$ws.Range("E7").Formula = "=""  -0.06%  """
$ws.Range("D8").Formula = "=""3.140.70"""
$ws.Range("E8").Formula = "=""  -4.78%  """
$ws.Range("E9").Formula = "=""  -4.54%  """
$ws.Range("E10").Formula = "=""  -7.56%  """
$ws.Range("E11").Formula = "=""  -8.10%  """
$ws.Range("E12").Formula = "=""  -8.17%  """
$ws.Range("D13").Formula = "=""3.683.57"""
$ws.Range("E13").Formula = "=""  -4.58%  """
$ws.Range("E14").Formula = "=""  -0.28%  """
$ws.Range("D15").Formula = "=""25.23"""
$ws.Range("E15").Formula = "=""  -5.16%  """
$ws.Range("D16").Formula = "=""3.148.19"""
$ws.Range("E16").Formula = "=""  -4.45%  """
$ws.Range("D17").Formula = "=""58.360.66"""
$ws.Range("E17").Formula = "=""  -3.12%  """
$ws.Range("E18").Formula = "=""  -7.15%  """
$ws.Range("E19").Formula = "=""  -4.87%  """
$ws.Range("D20").Formula = "=""13.01"""
$ws.Range("E20").Formula = "=""  -4.83%  """
$ws.Range("D21").Formula = "=""7.90"""
$ws.Range("E21").Formula = "=""  -6.99%  """
$ws.Range("D22").Formula = "=""342.69"""
$ws.Range("E22").Formula = "=""  -8.24%  """
$ws.Range("E23").Formula = "=""  +0.00%  """
$ws.Range("D24").Formula = "=""0.509"""
$ws.Range("E24").Formula = "=""  -3.98%  """
$ws.Range("D25").Formula = "=""67.51"""
$ws.Range("E25").Formula = "=""  -7.11%  """
$ws.Range("D26").Formula = "=""3.285.90"""
$ws.Range("E26").Formula = "=""  -4.19%  """
$ws.Range("E27").Formula = "=""  -1.35%  """
$ws.Range("E28").Formula = "=""  -0.22%  """
$ws.Range("D29").Formula = "=""0.0₃0937"""
$ws.Range("E29").Formula = "=""  -8.05%  """
$ws.Range("D30").Formula = "=""6.78"""
$ws.Range("E30").Formula = "=""  -3.13%  """
$ws.Range("E31").Formula = "=""  -0.10%  """
$ws.Range("E32").Formula = "=""  +4.58%  """
$ws.Range("E33").Formula = "=""  -7.34%  """
$ws.Range("D34").Formula = "=""6.88"""
$ws.Range("E34").Formula = "=""  -7.07%  """
$ws.Range("E35").Formula = "=""  -4.54%  """
$ws.Range("D36").Formula = "=""159.68"""
$ws.Range("E36").Formula = "=""  -3.82%  """
$ws.Range("E37").Formula = "=""  -3.90%  """
$ws.Range("D38").Formula = "=""6.21"""
$ws.Range("E38").Formula = "=""  -5.88%  """
$ws.Range("E39").Formula = "=""  -8.92%  """
$ws.Range("D40").Formula = "=""0.0686"""
$ws.Range("E40").Formula = "=""  -4.73%  """
$ws.Range("D41").Formula = "=""3.174.65"""
$ws.Range("E41").Formula = "=""  -4.55%  """
$ws.Range("D42").Formula = "=""40.41"""
$ws.Range("E42").Formula = "=""  -2.90%  """
$ws.Range("D43").Formula = "=""23.75"""
$ws.Range("E43").Formula = "=""  -6.20%  """
$ws.Range("D44").Formula = "=""0.693"""
$ws.Range("E44").Formula = "=""  -7.32%  """
$ws.Range("D45").Formula = "=""1.08"""
$ws.Range("E45").Formula = "=""  -2.64%  """
$ws.Range("D46").Formula = "=""3.93"""
$ws.Range("E46").Formula = "=""  -3.50%  """
$ws.Range("E47").Formula = "=""  -0.01%  """
$ws.Range("E48").Formula = "=""  -6.65%  """
$ws.Range("D49").Formula = "=""2.289.13"""
$ws.Range("E49").Formula = "=""  -1.27%  """
$ws.Range("D50").Formula = "=""6.17"""
$ws.Range("E50").Formula = "=""  -2.42%  """
$ws.Range("D51").Formula = "=""20.52"""
$ws.Range("E51").Formula = "=""  -4.37%  """

# Freeze all the temporary formulas above into plain text values in one pass.
$changedRange = $ws.Range("D2:E51")
$changedRange.Copy()
$changedRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

